$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: this is now an "Employees" sample dataset
$ws.Name = "Sample.Employees"

# --- Row 3: header row (propertyName, employee-1, employee-2) ---
$ws.Range("B3").Value = "propertyName"
$ws.Range("C3").Value = "employee-1"
$ws.Range("D3").Value = "employee-2"

# --- Row 4: employee_name ---
$ws.Range("B4").Value = "employee_name"
$ws.Range("C4").Value = "Kermet Frog"
$ws.Range("D4").Value = "Miss Piggy"

# --- Row 5: employee_salary (reuse the currency-style formatting from the
#            old price row so the new salary values keep the same look) ---
$ws.Range("C6:D6").Copy() | Out-Null
$ws.Range("C5:D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B5").Value = "employee_salary"
$ws.Range("C5").Value = 50000
$ws.Range("D5").Value = 45000
$ws.Range("E5").Clear()

# --- Row 6: employee_age ---
$ws.Range("B6").Value = "employee_age"
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 66

# --- Row 7: drop the old description row but keep the styled blank label cell ---
$ws.Range("C7").Clear()
$ws.Range("D7").Clear()
$ws.Range("B7").ClearContents()

# --- A few trailing blank rows, matching the expanded used range ---
$ws.Range("B9").Borders.LineStyle = -4142
$ws.Range("B10").Borders.LineStyle = -4142
$ws.Range("B11").Borders.LineStyle = -4142

# --- Column widths ---
$ws.Range("B:B").ColumnWidth = 21.67
$ws.Range("C:D").ColumnWidth = 16.33

# --- Selection ---
$ws.Range("C3:D6").Select()
